$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of effort data, reusing the date style already used by column A
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A19").Value = 41187
$ws.Range("B19").Value = 1.75
$ws.Range("D19").Value = "Code cleanup check of all test cases as preparation of re-implementation index->pointer "

# Move the selection like the original author did after entering the new row
$ws.Range("A20").Select()
